# "imported data into ipynb" -- clean up the State codes sheet:
#  - remove the duplicate "UNION TERRITORY" sub-header row
#  - unmerge the old B28:B29 merged cell (UTTARAKHAND / UK code + the
#    "(FORMERLY UTTARANCHAL)" note) so every row is a plain 2-column record
#  - sort the whole state/code list alphabetically by state name
#  - refresh column width / selection to match the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second "UNION TERRITORY" header row (was row 32) entirely.
$ws.Rows(32).Delete()

# The UTTARAKHAND code cell was merged with the row below it
# (the "(FORMERLY UTTARANCHAL)" note); split it back into individual cells
# and drop the left-alignment that only made sense while merged.
$ws.Range("B28:B29").UnMerge()
$ws.Range("B28:B29").HorizontalAlignment = 1

# Sort all the data (everything below the main header row) alphabetically
# by the state/territory name in column A.
$dataRange = $ws.Range("A2:B37")
$keyRange = $ws.Range("A2:A37")
$dataRange.Sort($keyRange, 1)

# Column A now holds longer entries ("DADRA AND NAGAR HAVELI", etc.) at
# visible rows, so widen it to keep the best-fit look.
$ws.Columns("A").ColumnWidth = 29.42

# The wrapped note row shrinks now that it sits under a wider column, and
# the header row's wrap height is refreshed too.
$ws.Rows(29).EntireRow.AutoFit()
$ws.Rows(1).RowHeight = 42
$ws.Rows(2).RowHeight = 28

# Reflect where the user ended up after the cleanup/sort.
$ws.Range("A6").Select()
